$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last 4 rows (model_9_0_21 .. model_9_0_24 no longer present)
$ws.Range("A23:I26").EntireRow.Delete() | Out-Null

$ws.Range("A2").Value2 = "model_9_0_0"
$ws.Range("B2").Value2 = -0.03972843564838335
$ws.Range("C2").Value2 = -0.0884385054446204
$ws.Range("D2").Value2 = -1.087144929704265
$ws.Range("E2").Value2 = -0.1642415882357695
$ws.Range("F2").Value2 = 1.150673031806946
$ws.Range("G2").Value2 = 1.13770854473114
$ws.Range("H2").Value2 = 0.2058674991130829
$ws.Range("I2").Value2 = 0.6991952061653137

$ws.Range("A3").Value2 = "model_9_0_1"
$ws.Range("B3").Value2 = 0.4916056504106544
$ws.Range("C3").Value2 = 0.3592414826205835
$ws.Range("D3").Value2 = -1.136531317309955
$ws.Range("E3").Value2 = 0.2444502059770303
$ws.Range("F3").Value2 = 0.5626426935195923
$ws.Range("G3").Value2 = 0.6697635054588318
$ws.Range("H3").Value2 = 0.2107387632131577
$ws.Range("I3").Value2 = 0.4537518322467804

$ws.Range("A4").Value2 = "model_9_0_2"
$ws.Range("B4").Value2 = 0.499150704253348
$ws.Range("C4").Value2 = 0.3660649567553653
$ws.Range("D4").Value2 = -1.07631835272905
$ws.Range("E4").Value2 = 0.2553914064614908
$ws.Range("F4").Value2 = 0.5542925596237183
$ws.Range("G4").Value2 = 0.6626311540603638
$ws.Range("H4").Value2 = 0.204799622297287
$ws.Range("I4").Value2 = 0.4471809864044189

$ws.Range("A5").Value2 = "model_9_0_3"
$ws.Range("B5").Value2 = 0.5773583985224029
$ws.Range("C5").Value2 = 0.4437560720335975
$ws.Range("D5").Value2 = -0.5521265489019009
$ws.Range("E5").Value2 = 0.3674935247027381
$ws.Range("F5").Value2 = 0.4677397310733795
$ws.Range("G5").Value2 = 0.5814232230186462
$ws.Range("H5").Value2 = 0.1530954539775848
$ws.Range("I5").Value2 = 0.3798571228981018

$ws.Range("A6").Value2 = "model_9_0_4"
$ws.Range("B6").Value2 = 0.6209935765155035
$ws.Range("C6").Value2 = 0.5197316178277978
$ws.Range("D6").Value2 = -0.5304422286118082
$ws.Range("E6").Value2 = 0.4391759340247908
$ws.Range("F6").Value2 = 0.4194484353065491
$ws.Range("G6").Value2 = 0.5020085573196411
$ws.Range("H6").Value2 = 0.1509566009044647
$ws.Range("I6").Value2 = 0.3368076384067535

$ws.Range("A7").Value2 = "model_9_0_5"
$ws.Range("B7").Value2 = 0.6372879415464328
$ws.Range("C7").Value2 = 0.5159088158002698
$ws.Range("D7").Value2 = -0.4418830706259695
$ws.Range("E7").Value2 = 0.4424980528110082
$ws.Range("F7").Value2 = 0.4014153480529785
$ws.Range("G7").Value2 = 0.5060043334960938
$ws.Range("H7").Value2 = 0.1422214955091476
$ws.Range("I7").Value2 = 0.3348125219345093

$ws.Range("A8").Value2 = "model_9_0_20"
$ws.Range("B8").Value2 = 0.6409602587791607
$ws.Range("C8").Value2 = 0.4177508671379512
$ws.Range("D8").Value2 = -3.111388339061044
$ws.Range("E8").Value2 = 0.1457267476257968
$ws.Range("F8").Value2 = 0.3973512053489685
$ws.Range("G8").Value2 = 0.6086056232452393
$ws.Range("H8").Value2 = 0.4055306613445282
$ws.Range("I8").Value2 = 0.5130410194396973

$ws.Range("A9").Value2 = "model_9_0_6"
$ws.Range("B9").Value2 = 0.6434646383184071
$ws.Range("C9").Value2 = 0.5195434904942479
$ws.Range("D9").Value2 = -0.6591275647140535
$ws.Range("E9").Value2 = 0.4290565676948064
$ws.Range("F9").Value2 = 0.3945796191692352
$ws.Range("G9").Value2 = 0.5022051334381104
$ws.Range("H9").Value2 = 0.1636496037244797
$ws.Range("I9").Value2 = 0.3428849279880524

$ws.Range("A10").Value2 = "model_9_0_19"
$ws.Range("B10").Value2 = 0.6439122794146472
$ws.Range("C10").Value2 = 0.4465993961604744
$ws.Range("D10").Value2 = -2.912075653107569
$ws.Range("E10").Value2 = 0.1877137103313175
$ws.Range("F10").Value2 = 0.3940841853618622
$ws.Range("G10").Value2 = 0.5784511566162109
$ws.Range("H10").Value2 = 0.3858712315559387
$ws.Range("I10").Value2 = 0.4878254532814026

$ws.Range("A11").Value2 = "model_9_0_8"
$ws.Range("B11").Value2 = 0.6481178669463624
$ws.Range("C11").Value2 = 0.5127434073537129
$ws.Range("D11").Value2 = -0.7476493134695508
$ws.Range("E11").Value2 = 0.415948943281196
$ws.Range("F11").Value2 = 0.3894298374652863
$ws.Range("G11").Value2 = 0.5093130469322205
$ws.Range("H11").Value2 = 0.1723810285329819
$ws.Range("I11").Value2 = 0.350756824016571

$ws.Range("A12").Value2 = "model_9_0_9"
$ws.Range("B12").Value2 = 0.6486544736355413
$ws.Range("C12").Value2 = 0.5118296422149399
$ws.Range("D12").Value2 = -0.7490322683219819
$ws.Range("E12").Value2 = 0.4150000634650024
$ws.Range("F12").Value2 = 0.3888359367847443
$ws.Range("G12").Value2 = 0.5102682113647461
$ws.Range("H12").Value2 = 0.1725174486637115
$ws.Range("I12").Value2 = 0.3513266742229462

$ws.Range("A13").Value2 = "model_9_0_7"
$ws.Range("B13").Value2 = 0.649101664982485
$ws.Range("C13").Value2 = 0.5133655307318852
$ws.Range("D13").Value2 = -0.6349814915156986
$ws.Range("E13").Value2 = 0.425230132291096
$ws.Range("F13").Value2 = 0.388341099023819
$ws.Range("G13").Value2 = 0.5086627602577209
$ws.Range("H13").Value2 = 0.1612679362297058
$ws.Range("I13").Value2 = 0.3451829254627228

$ws.Range("A14").Value2 = "model_9_0_10"
$ws.Range("B14").Value2 = 0.6499361848832776
$ws.Range("C14").Value2 = 0.5087765694462187
$ws.Range("D14").Value2 = -0.7323641600880564
$ws.Range("E14").Value2 = 0.4134748081062696
$ws.Range("F14").Value2 = 0.3874174952507019
$ws.Range("G14").Value2 = 0.5134594440460205
$ws.Range("H14").Value2 = 0.1708733737468719
$ws.Range("I14").Value2 = 0.3522426784038544

$ws.Range("A15").Value2 = "model_9_0_11"
$ws.Range("B15").Value2 = 0.6500040131789804
$ws.Range("C15").Value2 = 0.4963241057621021
$ws.Range("D15").Value2 = -0.7962722118566317
$ws.Range("E15").Value2 = 0.3970614867214067
$ws.Range("F15").Value2 = 0.3873424530029297
$ws.Range("G15").Value2 = 0.5264756083488464
$ws.Range("H15").Value2 = 0.1771769970655441
$ws.Range("I15").Value2 = 0.3620998561382294

$ws.Range("A16").Value2 = "model_9_0_12"
$ws.Range("B16").Value2 = 0.6506829329192172
$ws.Range("C16").Value2 = 0.4977609430705744
$ws.Range("D16").Value2 = -0.7772950870620619
$ws.Range("E16").Value2 = 0.3998523002398046
$ws.Range("F16").Value2 = 0.3865910172462463
$ws.Range("G16").Value2 = 0.5249737501144409
$ws.Range("H16").Value2 = 0.1753051728010178
$ws.Range("I16").Value2 = 0.3604238033294678

$ws.Range("A17").Value2 = "model_9_0_15"
$ws.Range("B17").Value2 = 0.6511930075805161
$ws.Range("C17").Value2 = 0.5000649807737921
$ws.Range("D17").Value2 = -0.9625664473823559
$ws.Range("E17").Value2 = 0.3876559037358623
$ws.Range("F17").Value2 = 0.3860265612602234
$ws.Range("G17").Value2 = 0.5225654244422913
$ws.Range("H17").Value2 = 0.1935795843601227
$ws.Range("I17").Value2 = 0.3677484393119812

$ws.Range("A18").Value2 = "model_9_0_14"
$ws.Range("B18").Value2 = 0.6515866023183426
$ws.Range("C18").Value2 = 0.5018230888387345
$ws.Range("D18").Value2 = -0.8022181248745075
$ws.Range("E18").Value2 = 0.4016691140632869
$ws.Range("F18").Value2 = 0.3855909705162048
$ws.Range("G18").Value2 = 0.5207276940345764
$ws.Range("H18").Value2 = 0.1777634620666504
$ws.Range("I18").Value2 = 0.3593326807022095

$ws.Range("A19").Value2 = "model_9_0_13"
$ws.Range("B19").Value2 = 0.652192054289948
$ws.Range("C19").Value2 = 0.500512597537667
$ws.Range("D19").Value2 = -0.7175411598861254
$ws.Range("E19").Value2 = 0.4070061697830383
$ws.Range("F19").Value2 = 0.3849209249019623
$ws.Range("G19").Value2 = 0.5220974683761597
$ws.Range("H19").Value2 = 0.1694112718105316
$ws.Range("I19").Value2 = 0.3561274707317352

$ws.Range("A20").Value2 = "model_9_0_16"
$ws.Range("B20").Value2 = 0.6596564326409486
$ws.Range("C20").Value2 = 0.4976843287336935
$ws.Range("D20").Value2 = -0.7922410294400484
$ws.Range("E20").Value2 = 0.3986265298675639
$ws.Range("F20").Value2 = 0.3766600787639618
$ws.Range("G20").Value2 = 0.5250537991523743
$ws.Range("H20").Value2 = 0.1767793744802475
$ws.Range("I20").Value2 = 0.3611599504947662

$ws.Range("A21").Value2 = "model_9_0_18"
$ws.Range("B21").Value2 = 0.6623172796094626
$ws.Range("C21").Value2 = 0.4934331081686266
$ws.Range("D21").Value2 = -1.705002377921071
$ws.Range("E21").Value2 = 0.3241620398673376
$ws.Range("F21").Value2 = 0.3737152516841888
$ws.Range("G21").Value2 = 0.5294974446296692
$ws.Range("H21").Value2 = 0.2668104469776154
$ws.Range("I21").Value2 = 0.4058802127838135

$ws.Range("A22").Value2 = "model_9_0_17"
$ws.Range("B22").Value2 = 0.6645767344346503
$ws.Range("C22").Value2 = 0.5101408537935153
$ws.Range("D22").Value2 = -1.072930676856797
$ws.Range("E22").Value2 = 0.3884100335507084
$ws.Range("F22").Value2 = 0.3712146878242493
$ws.Range("G22").Value2 = 0.5120334625244141
$ws.Range("H22").Value2 = 0.204465463757515
$ws.Range("I22").Value2 = 0.367295503616333
